$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.935.73'
$ws.Range("E2").Value = '  +12.45%  '
$ws.Range("D3").Value = '2.686.52'
$ws.Range("E3").Value = '  +14.72%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''518.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +9.08%  '
$ws.Range("D6").Value = '''161.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +11.00%  '
$ws.Range("D8").Value = '''0.615'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.53%  '
$ws.Range("D9").Value = '2.694.25'
$ws.Range("E9").Value = '  +14.77%  '
$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").Value = '''6.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +13.14%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '''0.107'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +11.46%  '
$ws.Range("D12").Value = '''0.352'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.58%  '
$ws.Range("E13").Value = '  +1.92%  '
$ws.Range("D14").Value = '3.164.07'
$ws.Range("E14").Value = '  +14.95%  '
$ws.Range("D15").Value = '61.496.63'
$ws.Range("E15").Value = '  +11.59%  '
$ws.Range("D16").Value = '''22.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +13.99%  '
$ws.Range("D17").Value = '''0.0000143'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +11.07%  '
$ws.Range("D18").Value = '2.693.62'
$ws.Range("E18").Value = '  +14.65%  '
$ws.Range("D19").Value = '''4.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.07%  '
$ws.Range("D20").Value = '''358.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +14.57%  '
$ws.Range("D21").Value = '''10.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +12.35%  '
$ws.Range("D22").Value = '''6.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +11.10%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").Value = '''61.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.74%  '
$ws.Range("E25").Value = '  +9.78%  '
$ws.Range("E26").Value = '  +11.96%  '
$ws.Range("D27").Value = '2.790.20'
$ws.Range("E27").Value = '  +13.87%  '
$ws.Range("D28").Value = '''1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").Value = '0.0₃0870'
$ws.Range("E29").Value = '  +17.82%  '
$ws.Range("D30").Value = '''7.65'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.84%  '
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("D32").Value = '''19.87'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.86%  '
$ws.Range("D33").Value = '''158.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.94%  '
$ws.Range("E34").Value = '  +9.42%  '
$ws.Range("E35").Value = '  +12.24%  '
$ws.Range("D36").Value = '''4.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +13.37%  '
$ws.Range("E37").Value = '  +12.92%  '
$ws.Range("D38").Value = '''0.880'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.52%  '
$ws.Range("D39").Value = '''1.51'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +15.18%  '
$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").Value = '''0.849'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +37.17%  '
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").Value = '''306.22'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +23.32%  '
$ws.Range("D42").Value = '''3.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.31%  '
$ws.Range("D43").Value = '''36.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.24%  '
$ws.Range("D44").Value = '''0.648'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +12.58%  '
$ws.Range("D45").Value = '''0.0588'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +14.42%  '
$ws.Range("E46").Value = '  +0.82%  '
$ws.Range("D47").Value = '''20.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +23.91%  '
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("D49").Value = '''5.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +15.39%  '
$ws.Range("D50").Value = '''0.0240'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.51%  '
$ws.Range("D51").Value = '2.053.73'
$ws.Range("E51").Value = '  +14.48%  '
